$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 3723.5
$ws.Range("I43").Value = 526.6667
$ws.Range("J43").Value = 4789.1113
$ws.Range("K43").Value = 526.6667
$ws.Range("L43").Value = 4789.1113
$ws.Range("M43").Value = -457.6667
$ws.Range("N43").Value = -4927.1113
$ws.Range("H88").Value = 3975
$ws.Range("I88").Value = 1966.6666
$ws.Range("J88").Value = 10000
$ws.Range("K88").Value = 1966.6666
$ws.Range("L88").Value = 10000
$ws.Range("M88").Value = -1560.6666
$ws.Range("N88").Value = -10812
$ws.Range("H91").Value = 3975
$ws.Range("I91").Value = 1966.6666
$ws.Range("J91").Value = 10000
$ws.Range("K91").Value = 1966.6666
$ws.Range("L91").Value = 10000
$ws.Range("M91").Value = -562.6666
$ws.Range("N91").Value = -12808
$ws.Range("H98").Value = 2925.1904
$ws.Range("I98").Value = 2950.205
$ws.Range("J98").Value = 2600
$ws.Range("K98").Value = 2950.205
$ws.Range("L98").Value = 2600
$ws.Range("M98").Value = -1452.205
$ws.Range("N98").Value = -5596
$ws.Range("H116").Value = 2570.95
$ws.Range("I116").Value = 2694.9375
$ws.Range("J116").Value = 2075
$ws.Range("K116").Value = 2694.9375
$ws.Range("L116").Value = 2075
$ws.Range("M116").Value = 747.0625
$ws.Range("N116").Value = -8959
$ws.Range("H122").Value = 2925.1904
$ws.Range("I122").Value = 2950.205
$ws.Range("J122").Value = 2600
$ws.Range("K122").Value = 8850.615
$ws.Range("L122").Value = 7800
$ws.Range("M122").Value = -6400.615
$ws.Range("N122").Value = -12700

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9445.313
$ws.Range("I32").Value = 6268.6
$ws.Range("K32").Value = 6268.6
$ws.Range("M32").Value = -5981.6
$ws.Range("H61").Value = 9525396
$ws.Range("I61").Value = 11112609
$ws.Range("K61").Value = 11112609
$ws.Range("M61").Value = -11112397
$ws.Range("H102").Value = 1969.7273
$ws.Range("I102").Value = 1901.125
$ws.Range("J102").Value = 2152.6667
$ws.Range("K102").Value = 1901.125
$ws.Range("L102").Value = 2152.6667
$ws.Range("M102").Value = -279.125
$ws.Range("N102").Value = -5396.6667
$ws.Range("H122").Value = 1314.24
$ws.Range("I122").Value = 1323.1111
$ws.Range("J122").Value = 1291.4286
$ws.Range("K122").Value = 3969.3333
$ws.Range("L122").Value = 3874.2858
$ws.Range("M122").Value = -1519.3333
$ws.Range("N122").Value = -8774.2858
$ws.Range("H136").Value = 9525396
$ws.Range("I136").Value = 11112609
$ws.Range("K136").Value = 33337827
$ws.Range("M136").Value = -33335277

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H122").Value = 55000
$ws.Range("J122").Value = 55000
$ws.Range("L122").Value = 55000
$ws.Range("N122").Value = -64800
$ws.Range("H134").Value = 47105.5
$ws.Range("I134").Value = 59143.05
$ws.Range("J134").Value = 1362.8
$ws.Range("K134").Value = 177429.15
$ws.Range("L134").Value = 4088.4
$ws.Range("M134").Value = -174894.15
$ws.Range("N134").Value = -9158.4

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H97").Value = 29999
$ws.Range("J97").Value = 29999
$ws.Range("L97").Value = 29999
$ws.Range("N97").Value = -31981
$ws.Range("H100").Value = 40390
$ws.Range("J100").Value = 40390
$ws.Range("L100").Value = 40390
$ws.Range("N100").Value = -42554
$ws.Range("H122").Value = 4808243
$ws.Range("I122").Value = 7353416
$ws.Range("J122").Value = 694.44446
$ws.Range("K122").Value = 22060248
$ws.Range("L122").Value = 2083.33338
$ws.Range("M122").Value = -22057798
$ws.Range("N122").Value = -6983.33338
$ws.Range("H132").Value = 1348.8055
$ws.Range("I132").Value = 1277.8387
$ws.Range("J132").Value = 1788.8
$ws.Range("K132").Value = 3833.5161
$ws.Range("L132").Value = 5366.4
$ws.Range("M132").Value = -1303.5161
$ws.Range("N132").Value = -10426.4

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 129.75
$ws.Range("I7").Value = 115.333336
$ws.Range("J7").Value = 138.4
$ws.Range("K7").Value = 346.000008
$ws.Range("L7").Value = 415.2
$ws.Range("M7").Value = -234.000008
$ws.Range("N7").Value = -639.2
$ws.Range("H23").Value = 173.66667
$ws.Range("I23").Value = 50
$ws.Range("J23").Value = 235.5
$ws.Range("K23").Value = 150
$ws.Range("L23").Value = 706.5
$ws.Range("M23").Value = 85
$ws.Range("N23").Value = -1176.5
$ws.Range("H131").Value = 2394.3967
$ws.Range("J131").Value = 1695.638
$ws.Range("L131").Value = 5086.914
$ws.Range("N131").Value = -15166.914
$ws.Range("H132").Value = 83334520
$ws.Range("I132").Value = 333333860
$ws.Range("J132").Value = 1398.7778
$ws.Range("K132").Value = 3000004740
$ws.Range("L132").Value = 12589.0002
$ws.Range("M132").Value = -3000002210
$ws.Range("N132").Value = -17649.0002

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H64").Value = 0
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 0
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H113").Value = 5689.577
$ws.Range("I113").Value = 9627.385
$ws.Range("J113").Value = 1751.7693
$ws.Range("K113").Value = 9627.385
$ws.Range("L113").Value = 1751.7693
$ws.Range("M113").Value = -7457.385
$ws.Range("N113").Value = -6091.7693
$ws.Range("H122").Value = 3154.422
$ws.Range("I122").Value = 2443.1035
$ws.Range("J122").Value = 4443.6875
$ws.Range("K122").Value = 7329.310500000001
$ws.Range("L122").Value = 13331.0625
$ws.Range("M122").Value = -4879.310500000001
$ws.Range("N122").Value = -18231.0625
$ws.Range("H127").Value = 54000
$ws.Range("J127").Value = 54000
$ws.Range("L127").Value = 54000
$ws.Range("N127").Value = -63920
$ws.Range("H132").Value = 2521.4285
$ws.Range("I132").Value = 2487.4
$ws.Range("J132").Value = 2606.5
$ws.Range("K132").Value = 7462.200000000001
$ws.Range("L132").Value = 7819.5
$ws.Range("M132").Value = -4932.200000000001
$ws.Range("N132").Value = -12879.5
$ws.Range("H136").Value = 15900.375
$ws.Range("J136").Value = 15900.375
$ws.Range("L136").Value = 47701.125
$ws.Range("N136").Value = -52801.125

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 1437.2106
$ws.Range("I61").Value = 1421.9286
$ws.Range("J61").Value = 1480
$ws.Range("K61").Value = 1421.9286
$ws.Range("L61").Value = 1480
$ws.Range("M61").Value = -1219.9286
$ws.Range("N61").Value = -1884
$ws.Range("H113").Value = 1437.2106
$ws.Range("I113").Value = 1421.9286
$ws.Range("J113").Value = 1480
$ws.Range("K113").Value = 1421.9286
$ws.Range("L113").Value = 1480
$ws.Range("M113").Value = 748.0714
$ws.Range("N113").Value = -5820
$ws.Range("H132").Value = 36632.97
$ws.Range("I132").Value = 40283.82
$ws.Range("J132").Value = 2558.3333
$ws.Range("K132").Value = 120851.46
$ws.Range("L132").Value = 7674.999899999999
$ws.Range("M132").Value = -118321.46
$ws.Range("N132").Value = -12734.9999

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H17").Value = 50000000
$ws.Range("I17").Value = 50000000
$ws.Range("K17").Value = 50000000
$ws.Range("M17").Value = -49999828
$ws.Range("H81").Value = 4406.8667
$ws.Range("I81").Value = 1243.8334
$ws.Range("J81").Value = 6515.5557
$ws.Range("K81").Value = 2487.6668
$ws.Range("L81").Value = 13031.1114
$ws.Range("M81").Value = -1426.6668
$ws.Range("N81").Value = -15153.1114
$ws.Range("H84").Value = 4406.8667
$ws.Range("I84").Value = 1243.8334
$ws.Range("J84").Value = 6515.5557
$ws.Range("K84").Value = 12438.334
$ws.Range("L84").Value = 65155.557
$ws.Range("M84").Value = -7134.333999999999
$ws.Range("N84").Value = -75763.557
$ws.Range("H132").Value = 11122.25
$ws.Range("I132").Value = 12842.923
$ws.Range("J132").Value = 3666
$ws.Range("K132").Value = 38528.769
$ws.Range("L132").Value = 10998
$ws.Range("M132").Value = -35998.769
$ws.Range("N132").Value = -16058
$ws.Range("H136").Value = 7189.4062
$ws.Range("I136").Value = 8460.615
$ws.Range("J136").Value = 1680.8334
$ws.Range("K136").Value = 25381.845
$ws.Range("L136").Value = 5042.5002
$ws.Range("M136").Value = -22831.845
$ws.Range("N136").Value = -10142.5002
